$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (row 2)
$ws.Range("C2").Value = 3.3815868999999998
$ws.Range("D2").Value = 4.9158626999999999
$ws.Range("E2").Value = 8.0134588000000004

# Update data values (row 3)
$ws.Range("C3").Value = 0.214587
$ws.Range("D3").Value = 0.204544
$ws.Range("E3").Value = 0.30090499999999998

# Update data values (row 4) -- written in plain decimal form to avoid
# scientific-notation literal parsing issues
$ws.Range("C4").Value = 0.042917499999999997
$ws.Range("D4").Value = 0.040908899999999998
$ws.Range("E4").Value = 0.080181000000000002

# Set explicit column width for column D so the saved sheet width matches
# column B's stored width of 10 (Excel's ColumnWidth property uses
# "characters" units that differ from the stored XML width units, so we
# request 9.17 characters to land on a stored width of 10)
$ws.Columns("D").ColumnWidth = 9.17

# Update the active cell selection
$ws.Range("E9").Select()
